# ---------------------------------------------------------------------------
# Adds a new "Player Info" sheet in front of the existing "ODI Batting" /
# "ODI Bowling" sheets, and replaces the MATCH_CARD_LINK column (full
# howstat.com URL) with a MATCH_CODE column (just the numeric match code)
# on both existing sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Player Info" worksheet as the first sheet in the book.
#    (Note: worksheet handles can become stale once the sheet collection is
#    mutated, so re-fetch sheets by name from $wb right before using them.)
# ---------------------------------------------------------------------------
$anchorSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($anchorSheet)
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $headers.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Re-use the same bold / bordered / centered header style already present
# in the workbook (style used by the other sheets' header row).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerData = @("6547", "Michael Gordeon Bracewell", "Left Handed", "Right Arm Off Break")
for ($c = 1; $c -le $playerData.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c)
    # Force text storage (rather than a numeric cell) for values that look
    # like numbers, then drop back to the default "Normal" style so no
    # stray number-format style is left applied to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $playerData[$c - 1]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace the full match-card
#    URL with just the numeric match code on both existing sheets.
# ---------------------------------------------------------------------------
function Update-MatchCodeColumn($sheet, $col) {
    $sheet.Cells.Item(1, $col).Value = "MATCH_CODE"

    $usedRows = $sheet.UsedRange.Rows.Count
    for ($r = 2; $r -le $usedRows; $r++) {
        $cell = $sheet.Cells.Item($r, $col)
        $val = $cell.Value()
        if ($val -match 'MatchCode=(\d+)') {
            $code = $matches[1]
            $cell.NumberFormat = "@"
            $cell.Value = $code
            $cell.Style = "Normal"
        }
    }
}

# "ODI Batting": MATCH_CARD_LINK is column D (4)
$battingSheet = $wb.Worksheets.Item("ODI Batting")
Update-MatchCodeColumn $battingSheet 4

# "ODI Bowling": MATCH_CARD_LINK is column B (2)
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
Update-MatchCodeColumn $bowlingSheet 2

Write-Host "Player Info sheet added; MATCH_CARD_LINK -> MATCH_CODE applied to ODI Batting and ODI Bowling."
